$wb = $excel.ActiveWorkbook

# --- Sheet 1: ROW50-FE-LIFTER -> add row 76 ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Cells.Item(76, 1).NumberFormat = $ws1.Cells.Item(75, 1).NumberFormat
$ws1.Cells.Item(76, 1).Value = 45762.26023787037
$ws1.Cells.Item(76, 2).Value = "0x01,0x90"
$ws1.Cells.Item(76, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"
$ws1.Cells.Item(76, 4).Value = "0x01,0x4e"
$ws1.Cells.Item(76, 5).Value = "0xe"
$ws1.Cells.Item(76, 6).Value = 400
$ws1.Cells.Item(76, 7).Value = [double]"5.68631262647114e+23"
$ws1.Cells.Item(76, 8).Value = 334
$ws1.Cells.Item(76, 9).Value = 14

# --- Sheet 2: ROW50-MID-LIFTER -> add row 78 (ID_DEC column kept as text) ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Cells.Item(78, 1).NumberFormat = $ws2.Cells.Item(77, 1).NumberFormat
$ws2.Cells.Item(78, 1).Value = 45762.22238425926
$ws2.Cells.Item(78, 2).Value = "0x01,0x90 "
$ws2.Cells.Item(78, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
$ws2.Cells.Item(78, 4).Value = "0x01,0x52"
$ws2.Cells.Item(78, 5).Value = "0x19"
$ws2.Cells.Item(78, 6).Value = 400
$ws2.Cells.Item(78, 7).NumberFormat = "@"
$ws2.Cells.Item(78, 7).Value = "568631262647113771663628"
$ws2.Cells.Item(78, 8).Value = 338
$ws2.Cells.Item(78, 9).Value = 25

# --- Sheet 3: ROW11-FE-LIFTER -> add row 76 ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Cells.Item(76, 1).NumberFormat = $ws3.Cells.Item(75, 1).NumberFormat
$ws3.Cells.Item(76, 1).Value = 45762.29122293981
$ws3.Cells.Item(76, 2).Value = "0x01,0x90"
$ws3.Cells.Item(76, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,"
$ws3.Cells.Item(76, 4).Value = "0x01,0x4e"
$ws3.Cells.Item(76, 5).Value = "0x14"
$ws3.Cells.Item(76, 6).Value = 400
$ws3.Cells.Item(76, 7).Value = [double]"5.68631262647114e+23"
$ws3.Cells.Item(76, 8).Value = 334
$ws3.Cells.Item(76, 9).Value = 20

# --- Sheet 4: ROW11-MID-LIFTER -> add row 76 ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Cells.Item(76, 1).NumberFormat = $ws4.Cells.Item(75, 1).NumberFormat
$ws4.Cells.Item(76, 1).Value = 45762.41875864584
$ws4.Cells.Item(76, 2).Value = "0x01,0x90"
$ws4.Cells.Item(76, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,"
$ws4.Cells.Item(76, 4).Value = "0x01,0x52"
$ws4.Cells.Item(76, 5).Value = "0x19"
$ws4.Cells.Item(76, 6).Value = 400
$ws4.Cells.Item(76, 7).Value = [double]"5.68631262647114e+23"
$ws4.Cells.Item(76, 8).Value = 338
$ws4.Cells.Item(76, 9).Value = 25
